$d = $word.ActiveDocument

$replacements = @(
    @("322÷3=", "703÷4="),
    @("215÷7=", "579÷8="),
    @("503÷3=", "248÷2="),
    @("826÷2=", "157÷9="),
    @("815÷3=", "119÷5="),
    @("608÷6=", "795÷3="),
    @("138÷3=", "951÷7="),
    @("796÷8=", "571÷8="),
    @("779÷3=", "795÷2="),
    @("995÷8=", "180÷9="),
    @("430÷4=", "688÷8="),
    @("508÷3=", "393÷9="),
    @("653÷6=", "303÷9="),
    @("952÷3=", "538÷6="),
    @("684÷9=", "226÷5="),
    @("976÷5=", "782÷2="),
    @("393÷7=", "410÷8="),
    @("521÷8=", "417÷6="),
    @("450÷8=", "609÷6="),
    @("823÷8=", "166÷8="),
    @("395÷3=", "782÷4="),
    @("561÷5=", "110÷7="),
    @("240÷9=", "549÷4="),
    @("429÷3=", "712÷5=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
